# Rows 2, 3 and 4 of the "Artfynd" sheet get cyclically rotated:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# Only the columns that actually carry per-record data differ between the
# three rows (the rest - C, I, P, S, T, U, V, W, AC, AD, AE, AG, AT, AW, AX,
# AY - are identical across rows 2-4 already), so only those columns need
# to move: A, B, D, E, F, G, H, Q, R, Y, AA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers for A, B, D, E, F, G, H, Q, R, Y, AA
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 25, 27)

# Use an out-of-range scratch row as temporary storage so the 3-way swap
# doesn't clobber data before it's been relocated. Using Range.Copy (rather
# than re-typing the literal value through .Value) preserves each cell's
# original type/formatting - e.g. date-like text such as "2023-07-04"
# stays plain text instead of being re-interpreted as a date serial.
$scratchRow = 1000

foreach ($col in $cols) {
  $ws.Cells.Item(2, $col).Copy($ws.Cells.Item($scratchRow, $col))
}
foreach ($col in $cols) {
  $ws.Cells.Item(3, $col).Copy($ws.Cells.Item(2, $col))
}
foreach ($col in $cols) {
  $ws.Cells.Item(4, $col).Copy($ws.Cells.Item(3, $col))
}
foreach ($col in $cols) {
  $ws.Cells.Item($scratchRow, $col).Copy($ws.Cells.Item(4, $col))
}
foreach ($col in $cols) {
  $ws.Cells.Item($scratchRow, $col).Clear()
}

Write-Host "Rotated rows 2-4 (A,B,D,E,F,G,H,Q,R,Y,AA)"
